# BOM updated to Rev 0.4
# Remove four obsolete BOM line items (C15/4.7uF cap, L1/2.2uH inductor,
# R7/453k resistor, U4/AP3417C regulator) from the SAMPad BOM table, and
# fold the references that pointed at the removed parts (R8 and C16) out
# of the two remaining grouped-reference rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-up so earlier row numbers stay valid while we work.
$ws.Rows(28).Delete()   # U4, AP3417C, IC REG BUCK ADJ 1A SYNC SOT25
$ws.Rows(21).Delete()   # R7, 453k, RES SMD 453K OHM 1% 1/8W 0805
$ws.Rows(10).Delete()   # L1, 2.2uH, FIXED IND 2.2UH 2.9A 44 MOHM SMD
$ws.Rows(4).Delete()    # C15, 4.7uF, CAP CER 4.7UF 6.3V X5R 0805

# Row 5 (C9, C16 - C18 / Qty 4) shifted up to row 4 after the C15 delete;
# drop the now-removed C16 reference and update the quantity to 3.
$ws.Range("A4").Value = "C9, C17 - C18"
$ws.Range("C4").Value = 3

# Row 22 (R8, R9, R14 / Qty 3) shifted up to row 19 after the deletes;
# drop the now-removed R8 reference and update the quantity to 2.
$ws.Range("A19").Value = "R9, R14"
$ws.Range("C19").Value = 2

# Update the hidden AutoFilter defined name so it reflects the new, smaller
# table extent (rows 1-26 of the now-shorter SAMPad sheet) instead of the
# stale A1:F30.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "SAMPad!_FilterDatabase") {
        $n.RefersTo = "=SAMPad!`$A`$1:`$F`$26"
    }
}

# Leave the cursor/selection on C5, matching the saved view state.
$ws.Range("C5").Select()
